{"js": "// Edit: update \"Multispectral Sensor comparison\" table.\n//  1. Spatial resolution / PlanetScope cell: \"3-4 ,\" -> two runs \"3-4\" + \" m\"\n//  2. Temporal resolution row: fill in the four empty data cells with revisit\n//     cadence text for each sensor.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells + the text of the first cell (used to identify the\n// row by its label, independent of row ordering).\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet spatialResolutionRow = null;\nlet temporalResolutionRow = null;\nfor (const row of rows.items) {\n  const labelCell = row.cells.items[0];\n  const label = labelCell.body.text;\n  if (label.indexOf(\"Spatial resolution\") === 0) {\n    spatialResolutionRow = row;\n  } else if (label.indexOf(\"Temporal resolution\") === 0) {\n    temporalResolutionRow = row;\n  }\n}\n\n// 1) Split \"3-4 ,\" into \"3-4\" and \" m\" as two separate runs in the\n//    PlanetScope column (3rd cell, index 2) of the Spatial resolution row.\nconst planetScopeCell = spatialResolutionRow.cells.items[2];\nconst splitOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r><w:t>3-4</w:t></w:r><w:r><w:t xml:space=\"preserve\"> m</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nplanetScopeCell.body.insertOoxml(splitOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Fill in the four empty \"Temporal resolution\" data cells.\nconst revisitText = [\n  \"Every 5 days or so\",\n  \"Near-daily (depends on cloud cover, depend on latitude)\",\n  \"As often as you fly the drones \",\n  \"Daily, hourly\",\n];\nconst temporalCells = temporalResolutionRow.cells.items;\nfor (let i = 0; i < revisitText.length; i++) {\n  temporalCells[i + 1].body.insertText(revisitText[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Edit: update \"Multispectral Sensor comparison\" table.\n#  1. Spatial resolution / PlanetScope cell: \"3-4 ,\" -> two runs \"3-4\" + \" m\"\n#  2. Temporal resolution row: fill in the four empty data cells with revisit\n#     cadence text for each sensor.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$spatialResolutionRow = $null\n$temporalResolutionRow = $null\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  $label = $t.Cell($r, 1).Range.Text\n  if ($label.StartsWith(\"Spatial resolution\")) {\n    $spatialResolutionRow = $r\n  } elseif ($label.StartsWith(\"Temporal resolution\")) {\n    $temporalResolutionRow = $r\n  }\n}\n\n# 1) Split \"3-4 ,\" into \"3-4\" and \" m\" as two separate runs in the\n#    PlanetScope column (3rd column) of the Spatial resolution row.\n$planetScopeCell = $t.Cell($spatialResolutionRow, 3)\n$cellRange = $planetScopeCell.Range\n$cellRange.MoveEnd(1, -1) | Out-Null\n$cellRange.Text = \"\"\n$cellRange.InsertXML('<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>3-4</w:t></w:r><w:r><w:t xml:space=\"preserve\"> m</w:t></w:r></w:p>')\n\n# 2) Fill in the four empty \"Temporal resolution\" data cells.\n$revisitText = @(\n  \"Every 5 days or so\",\n  \"Near-daily (depends on cloud cover, depend on latitude)\",\n  \"As often as you fly the drones \",\n  \"Daily, hourly\"\n)\nfor ($i = 0; $i -lt $revisitText.Length; $i++) {\n  $cell = $t.Cell($temporalResolutionRow, $i + 2)\n  $r = $cell.Range\n  $r.MoveEnd(1, -1) | Out-Null\n  $r.Text = $revisitText[$i]\n}\n"}
